$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100 (shifts existing rows 100.. down by one,
# matching the diff where old row 100 becomes new row 101, etc., and a
# brand-new record appears at row 100).
$ws.Rows("100:100").Insert()

# Populate the newly inserted row 100 with the new weekly record.
$ws.Range("A100").Value = 11
$ws.Range("B100").Value = "Vega Monumental Concepción"
$ws.Range("C100").Value = "Bíobío"
$ws.Range("D100").Value = 44894
$ws.Range("E100").Value = 8
$ws.Range("F100").Value = 100112003
$ws.Range("G100").Value = "Ajo"
$ws.Range("H100").Value = "Chino"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 14000
$ws.Range("L100").Value = 15000
$ws.Range("M100").Value = 14600
$ws.Range("N100").Value = "$/caja 10 kilos"
$ws.Range("O100").Value = "China"
$ws.Range("P100").Value = 1460
$ws.Range("Q100").Value = 10
$ws.Range("R100").Value = "Hortaliza"
